$wb = $excel.ActiveWorkbook

# --- Add "test_Smoke_Inspections_Superint" as a copy of "test_Smoke_Inspections_Admin",
#     placed at the end of the workbook ---
$srcInspections = $wb.Worksheets.Item("test_Smoke_Inspections_Admin")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcInspections.Copy($null, $lastSheet)
$newInspections = $wb.Worksheets.Item($wb.Worksheets.Count)
$newInspections.Name = "test_Smoke_Inspections_Superint"

# --- Add "test_Smoke_WorkOrders_Superinte" as a copy of "test_Smoke_WorkOrders_Admin",
#     placed at the end of the workbook ---
$srcWorkOrders = $wb.Worksheets.Item("test_Smoke_WorkOrders_Admin")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcWorkOrders.Copy($null, $lastSheet2)
$newWorkOrders = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWorkOrders.Name = "test_Smoke_WorkOrders_Superinte"

# The newly added sheet is now the active / selected tab.
# Update its selected cell to match the saved workbook state.
$newWorkOrders.Activate()
$newWorkOrders.Range("B8").Select()
